$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric/percentage-looking cells (Price, Volume%, Hora) ---
# These must be stored as literal text (trailing zeros / "%" suffix / "--" placeholders
# are significant), so force the cell to Text format before assigning the string value.
$textCells = @(
    @{ Ref = 'D2'; Value = '302.67' }
    @{ Ref = 'E2'; Value = '1.95%' }
    @{ Ref = 'G2'; Value = '16' }
    @{ Ref = 'D3'; Value = '31.90' }
    @{ Ref = 'E3'; Value = '1.68%' }
    @{ Ref = 'G3'; Value = '16' }
    @{ Ref = 'E4'; Value = '0.56%' }
    @{ Ref = 'G4'; Value = '16' }
    @{ Ref = 'D5'; Value = '0.07809' }
    @{ Ref = 'E5'; Value = '-2.59%' }
    @{ Ref = 'G5'; Value = '16' }
    @{ Ref = 'D6'; Value = '2.257' }
    @{ Ref = 'E6'; Value = '-18.58%' }
    @{ Ref = 'G6'; Value = '16' }
    @{ Ref = 'D7'; Value = '7.830' }
    @{ Ref = 'E7'; Value = '0.43%' }
    @{ Ref = 'G7'; Value = '16' }
    @{ Ref = 'E8'; Value = '0.52%' }
    @{ Ref = 'G8'; Value = '16' }
    @{ Ref = 'D9'; Value = '0.9234' }
    @{ Ref = 'E9'; Value = '-0.35%' }
    @{ Ref = 'G9'; Value = '16' }
    @{ Ref = 'D10'; Value = '0.1760' }
    @{ Ref = 'E10'; Value = '1.01%' }
    @{ Ref = 'G10'; Value = '16' }
    @{ Ref = 'D11'; Value = '0.07747' }
    @{ Ref = 'E11'; Value = '7.25%' }
    @{ Ref = 'G11'; Value = '16' }
    @{ Ref = 'D12'; Value = '0.08862' }
    @{ Ref = 'E12'; Value = '-1.15%' }
    @{ Ref = 'G12'; Value = '16' }
    @{ Ref = 'D13'; Value = '0.03148' }
    @{ Ref = 'E13'; Value = '4.98%' }
    @{ Ref = 'G13'; Value = '16' }
    @{ Ref = 'D14'; Value = '0.1000' }
    @{ Ref = 'E14'; Value = '-0.10%' }
    @{ Ref = 'G14'; Value = '16' }
    @{ Ref = 'D15'; Value = '0.001516' }
    @{ Ref = 'E15'; Value = '1.30%' }
    @{ Ref = 'G15'; Value = '16' }
    @{ Ref = 'D16'; Value = '0.005915' }
    @{ Ref = 'E16'; Value = '-0.67%' }
    @{ Ref = 'G16'; Value = '16' }
    @{ Ref = 'D17'; Value = '3.440' }
    @{ Ref = 'E17'; Value = '-2.77%' }
    @{ Ref = 'G17'; Value = '16' }
    @{ Ref = 'D18'; Value = '2.254' }
    @{ Ref = 'E18'; Value = '0.29%' }
    @{ Ref = 'G18'; Value = '16' }
    @{ Ref = 'G19'; Value = '16' }
    @{ Ref = 'D20'; Value = '0.1329' }
    @{ Ref = 'E20'; Value = '-1.20%' }
    @{ Ref = 'G20'; Value = '16' }
    @{ Ref = 'D21'; Value = '4.263' }
    @{ Ref = 'E21'; Value = '7.45%' }
    @{ Ref = 'G21'; Value = '16' }
    @{ Ref = 'D22'; Value = '0.1819' }
    @{ Ref = 'E22'; Value = '12.05%' }
    @{ Ref = 'G22'; Value = '16' }
    @{ Ref = 'D23'; Value = '0.04590' }
    @{ Ref = 'E23'; Value = '0.02%' }
    @{ Ref = 'G23'; Value = '16' }
    @{ Ref = 'D24'; Value = '0.001250' }
    @{ Ref = 'E24'; Value = '0.66%' }
    @{ Ref = 'G24'; Value = '16' }
    @{ Ref = 'D25'; Value = '0.004478' }
    @{ Ref = 'E25'; Value = '1.41%' }
    @{ Ref = 'G25'; Value = '16' }
    @{ Ref = 'D26'; Value = '0.0001251' }
    @{ Ref = 'E26'; Value = '4.42%' }
    @{ Ref = 'G26'; Value = '16' }
    @{ Ref = 'D27'; Value = '--' }
    @{ Ref = 'E27'; Value = '--%' }
    @{ Ref = 'G27'; Value = '16' }
    @{ Ref = 'G28'; Value = '16' }
    @{ Ref = 'G29'; Value = '16' }
    @{ Ref = 'G30'; Value = '16' }
    @{ Ref = 'G31'; Value = '16' }
    @{ Ref = 'G32'; Value = '16' }
    @{ Ref = 'G33'; Value = '16' }
    @{ Ref = 'G34'; Value = '16' }
    @{ Ref = 'G35'; Value = '16' }
    @{ Ref = 'G36'; Value = '16' }
    @{ Ref = 'G37'; Value = '16' }
    @{ Ref = 'G38'; Value = '16' }
    @{ Ref = 'D39'; Value = '0.01769' }
    @{ Ref = 'E39'; Value = '0.25%' }
    @{ Ref = 'G39'; Value = '16' }
    @{ Ref = 'D40'; Value = '0.04810' }
    @{ Ref = 'E40'; Value = '7.09%' }
    @{ Ref = 'G40'; Value = '16' }
    @{ Ref = 'D41'; Value = '0.007185' }
    @{ Ref = 'E41'; Value = '5.05%' }
    @{ Ref = 'G41'; Value = '16' }
    @{ Ref = 'E42'; Value = '1.56%' }
    @{ Ref = 'G42'; Value = '16' }
    @{ Ref = 'D43'; Value = '0.002122' }
    @{ Ref = 'E43'; Value = '-3.37%' }
    @{ Ref = 'G43'; Value = '16' }
    @{ Ref = 'D44'; Value = '0.009926' }
    @{ Ref = 'E44'; Value = '0.78%' }
    @{ Ref = 'G44'; Value = '16' }
    @{ Ref = 'D45'; Value = '0.00006264' }
    @{ Ref = 'E45'; Value = '-4.06%' }
    @{ Ref = 'G45'; Value = '16' }
    @{ Ref = 'E46'; Value = '0.15%' }
    @{ Ref = 'G46'; Value = '16' }
    @{ Ref = 'D47'; Value = '0.003567' }
    @{ Ref = 'E47'; Value = '-59.20%' }
    @{ Ref = 'G47'; Value = '16' }
    @{ Ref = 'D48'; Value = '1.171' }
    @{ Ref = 'E48'; Value = '42.68%' }
    @{ Ref = 'G48'; Value = '16' }
    @{ Ref = 'D49'; Value = '0.00002102' }
    @{ Ref = 'E49'; Value = '0.15%' }
    @{ Ref = 'G49'; Value = '16' }
    @{ Ref = 'D50'; Value = '0.0002002' }
    @{ Ref = 'E50'; Value = '0.15%' }
    @{ Ref = 'G50'; Value = '16' }
    @{ Ref = 'G51'; Value = '16' }
)
foreach ($item in $textCells) {
    $rng = $ws.Range($item.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# --- Plain text cells (Coin name, Link) ---
$plainCells = @(
    @{ Ref = 'B27'; Value = 'Spectre.aiUtilityToken' }
    @{ Ref = 'C27'; Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut' }
    @{ Ref = 'B28'; Value = 'LegolasExchange' }
    @{ Ref = 'C28'; Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo' }
    @{ Ref = 'B29'; Value = 'BitZToken' }
    @{ Ref = 'C29'; Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz' }
    @{ Ref = 'B30'; Value = 'Birake' }
    @{ Ref = 'C30'; Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir' }
    @{ Ref = 'B31'; Value = 'NashExchange' }
    @{ Ref = 'C31'; Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex' }
    @{ Ref = 'B32'; Value = 'AAXToken' }
    @{ Ref = 'C32'; Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab' }
    @{ Ref = 'B33'; Value = 'CenX' }
    @{ Ref = 'C33'; Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx' }
    @{ Ref = 'B34'; Value = 'BNIXToken' }
    @{ Ref = 'C34'; Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix' }
    @{ Ref = 'B35'; Value = 'UpBots' }
    @{ Ref = 'C35'; Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt' }
)
foreach ($item in $plainCells) {
    $ws.Range($item.Ref).Value = $item.Value
}

Write-Output "Applied $($textCells.Count) text-forced cell updates and $($plainCells.Count) plain cell updates."